$wb = $excel.ActiveWorkbook

# --- Make "Test Steps" the active sheet (was "Test Cases") ---
$ws1 = $wb.Worksheets.Item("Test Steps")
$ws1.Activate()

# --- Insert a new section header row above the existing data (row 2) ---
$ws1.Rows.Item(2).Insert() | Out-Null

$ws1.Range("A2").Value = "Web Application Testing"
$ws1.Range("A2:F2").Font.Bold = $true
$ws1.Range("A2:F2").Borders.LineStyle = 1

# --- Fix up the two mailto hyperlinks, which shifted down one row ---
$ws1.Range("A1").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("F6"), "mailto:ravikaanthe@rediffmail.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("F8"), "mailto:test@123") | Out-Null

# --- Selection ends up on B2 after the edit ---
$ws1.Range("B2").Select() | Out-Null
